$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression - only B2 changes slightly
$ws.Range("B2").Value = 14252680789404700

# Row 3: RandomForestRegressor - values change
$ws.Range("B3").Value = 396653422634.5623
$ws.Range("C3").Value = 356654758167.214
$ws.Range("D3").Value = 1775170097897886

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 458318030355.061
$ws.Range("C4").Value = 379616954435.5107
$ws.Range("D4").Value = 458318030355.0984

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 57565325496205.09
$ws.Range("C5").Value = 76649235372063.02
$ws.Range("D5").Value = 433340205936479.9
